$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Title (appears twice: H1 and bold paragraph near end)
while ($d.Content.Find.Execute("Play Medieval Money Slot Free - Review & Bonus Features", $true, $false, $false, $false, $false, $true, 1, $false, "Play Medieval Money Slot Game for Free", 2)) {}

# "What we like" bullet list
Replace-Text "Medieval themed graphics are appealing" "Medieval theme adds to the immersive gameplay experience"
Replace-Text "Offers 5 different bonus games" "Multiple bonus features offer exciting opportunities to win"
Replace-Text "Great winning potential, especially during the Joust Free Spins Bonus" "Wide bet range accommodates players with different budgets"
Replace-Text "Accessible to players with varying budgets" "High winning potential, especially during the Joust Free Spins Bonus"

# "What we don't like" bullet list
Replace-Text "Only available on desktop devices" "Cartoonish graphics may not appeal to all players"
Replace-Text "Graphics could be better" "Game only available on desktop devices"

# Meta description (italic paragraph at the end)
Replace-Text "Find out how to play Medieval Money, a medieval-themed slot game developed by IGT. Play it for free and enjoy 5 bonus features and great winning potential." "Read our review of Medieval Money, a slot game with a medieval theme and multiple bonus features. Play now for free."
